# Table S1 revision: add the "No ray bending" and "ray bending" location
# blocks (coordinates / TAT / water velocity) around the pre-existing
# bookmark paragraph.
#
# Word merges consecutive same-formatted runs that are typed/inserted back
# to back into a single run, so to reproduce the source document's distinct
# <w:r> boundaries for multi-run lines we build each run in its own
# temporary paragraph and then delete the paragraph mark that separates
# them -- joining the paragraphs back together while keeping each piece of
# text in its own run, exactly like Word does when you delete a pilcrow
# between two paragraphs.

$d = $word.ActiveDocument

function Insert-ParagraphAt([int]$Index, [string[]]$Parts) {
    # Inserts a new paragraph built out of one-or-more runs so that it
    # becomes paragraph number $Index (the paragraph currently occupying
    # that slot, and everything after it, shifts down).
    $anchor = $d.Paragraphs($Index).Range
    $anchor.Collapse(1)

    for ($i = 0; $i -lt $Parts.Length; $i++) {
        $anchor.InsertParagraphBefore()
    }

    for ($i = 0; $i -lt $Parts.Length; $i++) {
        $d.Paragraphs($Index + $i).Range.InsertAfter($Parts[$i])
    }

    # Re-join the scratch paragraphs into one, run-by-run, by deleting the
    # paragraph mark between each consecutive pair.
    for ($i = 0; $i -lt ($Parts.Length - 1); $i++) {
        $p = $d.Paragraphs($Index)
        $pEnd = $p.Range.End
        $mark = $d.Range($pEnd - 1, $pEnd)
        $mark.Delete()
    }
}

function Insert-EmptyParagraphAt([int]$Index) {
    $anchor = $d.Paragraphs($Index).Range
    $anchor.Collapse(1)
    $anchor.InsertParagraphBefore()
}

# Paragraph index of the original (only) paragraph -- the one carrying the
# _GoBack bookmark. Everything below is inserted relative to it; the index
# is advanced as content is added in front of it.
$bookmarkIndex = 1

Insert-ParagraphAt $bookmarkIndex @("EC03 ", "No", " ", "ray", " ", "bending")
$bookmarkIndex++

Insert-ParagraphAt $bookmarkIndex @("X:     -291.2 m (1.4) ")
$bookmarkIndex++

Insert-ParagraphAt $bookmarkIndex @("Y:    -170.5 m (2.5) ")
$bookmarkIndex++

Insert-ParagraphAt $bookmarkIndex @("Depth: -4742.0 m (5.4) ")
$bookmarkIndex++

Insert-ParagraphAt $bookmarkIndex @("TAT:   14.0 ", "ms", " (0.000000) ")
$bookmarkIndex++

Insert-ParagraphAt $bookmarkIndex @("Water Vel.: 1506.4 m/s (1.605341)")
$bookmarkIndex++

Insert-EmptyParagraphAt $bookmarkIndex
$bookmarkIndex++

Insert-ParagraphAt $bookmarkIndex @("EC03 ", "ray bending")
$bookmarkIndex++

# The bookmark paragraph itself gains a leading run of text, inserted right
# before the existing bookmarkStart/bookmarkEnd pair (still the same
# paragraph -- this is the paragraph shown in the diff context). InsertBefore
# (rather than InsertAfter) is required here so the new run lands in front of
# the bookmark markers instead of behind them.
$bookmarkRange = $d.Paragraphs($bookmarkIndex).Range
$bookmarkRange.Collapse(1)
$bookmarkRange.InsertBefore("X:     -291.3 m (1.5) ")
$bookmarkIndex++

Insert-ParagraphAt $bookmarkIndex @("Y:    -170.4 m (2.6) ")
$bookmarkIndex++

Insert-ParagraphAt $bookmarkIndex @("Depth: -4742.4 m (5.5) ")
$bookmarkIndex++

Insert-ParagraphAt $bookmarkIndex @("TAT:   14.0 ", "ms", " (0.000000) ")
$bookmarkIndex++

Insert-ParagraphAt $bookmarkIndex @("Water Vel.: 1506.6 m/s (1.638528)")
$bookmarkIndex++

Insert-EmptyParagraphAt $bookmarkIndex
$bookmarkIndex++

Insert-EmptyParagraphAt $bookmarkIndex
$bookmarkIndex++
